$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "Опыт р" + "а" + "боты" (three runs, identical rPr) -> a
# single run "Опыт работы". Searching for the exact concatenation of
# the three runs' text and "replacing" it with itself makes the
# engine collapse the matched range into a single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Опыт работы", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Опыт работы", 2) | Out-Null

# ------------------------------------------------------------------
# Change 2: "Асессор. Яндекс. С 02.11.2022 " + "по сей день." (two
# runs) -> "Асессор. Яндекс. " + "Ноябрь 2022 – Январь 2023" (two
# runs; same boundary/formatting as before, first run keeps its
# rsidRPr attribute).
# ------------------------------------------------------------------

# Locate (read-only, Replace:=wdReplaceNone) the end of the text that
# must stay untouched ("Асессор. Яндекс. ") and the bounds of the
# trailing sentinel text ("по сей день.") so the exact amount of text
# between them ("С 02.11.2022 ") is known without hard-coded offsets.
$prefixRange = $d.Content
$prefixRange.Find.Execute("Асессор. Яндекс. ", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0) | Out-Null
$prefixEnd = $prefixRange.End

$tailRange = $d.Content
$tailRange.Find.Execute("по сей день.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$tailStart = $tailRange.Start

# Replace the trailing run's text only; this run sits entirely after
# the text that must be preserved, so nothing in the first run is
# ever touched by this call. When the two runs are coalesced here,
# the merged run keeps the first run's rsidRPr attribute.
$tailRange.Text = "Ноябрь 2022 – Январь 2023"

# The whole sentence now lives in a single run. Re-split it exactly
# at the boundary we need ("Асессор. Яндекс. " | rest) by toggling a
# character formatting property on the trailing part - this creates
# a genuine run split without re-coalescing anything afterwards.
$para = $d.Paragraphs.Item($d.Paragraphs.Count)
$splitPoint = $prefixEnd
$rightPart = $d.Range($splitPoint, $para.Range.End - 1)
$rightPart.Font.Bold = 1

# Delete the now-isolated middle text ("С 02.11.2022 ") that sits
# between the preserved prefix and the new date text, without
# touching the (now separate) first run.
$midLen = $tailStart - $prefixEnd
if ($midLen -gt 0) {
    $midRange = $d.Range($splitPoint, $splitPoint + $midLen)
    $midRange.Text = ""
}

# Restore the formatting on the trailing run so it matches the rest
# of the paragraph's direct formatting again; the run boundary
# created above remains in place.
$para2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rightPart2 = $d.Range($splitPoint, $para2.Range.End - 1)
$rightPart2.Font.Bold = 0

Write-Output "Done"
